$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeployNginx")

$nginxConf = "server {" + "`n" + `
"    listen 80;" + "`n" + `
"    location / {" + "`n" + `
"        proxy_pass http://localhost:5000;" + "`n" + `
"        proxy_http_version 1.1;" + "`n" + `
"        proxy_set_header Upgrade `$http_upgrade;" + "`n" + `
"        proxy_set_header Connection keep-alive;" + "`n" + `
"        proxy_set_header Host `$host;" + "`n" + `
"        proxy_cache_bypass `$http_upgrade;" + "`n" + `
"    }" + "`n" + `
"    location /api/chat {" + "`n" + `
"        proxy_pass http://localhost:5000;" + "`n" + `
"        proxy_http_version 1.1;" + "`n" + `
"        proxy_set_header Upgrade `$http_upgrade;" + "`n" + `
"        proxy_set_header Connection `"upgrade`";" + "`n" + `
"        proxy_set_header Host `$host;" + "`n" + `
"        proxy_cache_bypass `$http_upgrade;" + "`n" + `
"    }" + "`n" + `
"}"

$description = "Thay đổi location /api/chat mapHub." + "`n" + "Sử dụng proxy_set_header Connection `"upgrade`"; with websocket realtime."

$ws.Range("A10").Value = "Deploy with signalr nginx"
$ws.Range("C10").Value = $nginxConf
$ws.Range("E10").Value = "https://medium.com/@alm.ozdmr/deployment-of-signalr-with-nginx-daf392cf2b93"
$ws.Range("B10").Value = $description

$ws.Hyperlinks.Add($ws.Range("E10"), "https://medium.com/@alm.ozdmr/deployment-of-signalr-with-nginx-daf392cf2b93")

$ws.Range("A10:D15").Font.Size = 14
$ws.Range("A10:D15").Font.Name = "Arial"
$ws.Range("A10:D15").WrapText = $true
$ws.Range("A10:D15").VerticalAlignment = -4108
$ws.Range("A10:D15").HorizontalAlignment = -4131

$ws.Range("E11:E15").Font.Size = 14
$ws.Range("E11:E15").Font.Name = "Arial"
$ws.Range("E11:E15").WrapText = $true
$ws.Range("E11:E15").VerticalAlignment = -4108
$ws.Range("E11:E15").HorizontalAlignment = -4131

$ws.Range("A10:E10").RowHeight = 332.5

$ws.ListObjects.Item(1).Resize($ws.Range("A2:E15"))

$ws.Range("A10").Select()
